$wb = $excel.ActiveWorkbook

# "Metadata" sheet holds the Property/Value table; set the "Name" row's value
$ws1 = $wb.Worksheets.Item("Metadata")
$ws1.Range("B4").Value = "CompetencespecifiqueVs"

# Update the "Date" row's value to reflect the new generation timestamp
$ws1.Range("B8").Value = "2025-07-18T06:40:38+00:00"
